$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers (new shared strings / text values)
$ws.Range("L1").Value = "hzj-混合调节_20170516_152754_ASIC_EEG"
$ws.Range("M1").Value = "hzj-混合调节_20170518_134207_ASIC_EEG"
$ws.Range("N1").Value = "hzj-混合调节_20170519_135415_ASIC_EEG"
$ws.Range("O1").Value = "zyx-混合调节_20170516_111228_ASIC_EEG"
$ws.Range("P1").Value = "zyx-混合调节_20170517_110944_ASIC_EEG"
$ws.Range("Q1").Value = "zyx-混合调节_20170518_112337_ASIC_EEG"
$ws.Range("R1").Value = "zyx-混合调节_20170519_124954_ASIC_EEG"
$ws.Range("S1").Value = "zyx-混合调节_20170522_111557_ASIC_EEG"

# Row 2 numeric values
$ws.Range("L2").Value = 0.64236111111111116
$ws.Range("M2").Value = 0.67697594501718217
$ws.Range("N2").Value = 0.71535580524344566
$ws.Range("O2").Value = 0.67628205128205132
$ws.Range("P2").Value = 0.76602564102564097
$ws.Range("Q2").Value = 0.74757281553398058
$ws.Range("R2").Value = 0.65161290322580645
$ws.Range("S2").Value = 0.7063106796116505

# Row 3 numeric values
$ws.Range("L3").Value = 0.62637362637362637
$ws.Range("M3").Value = 0.6387096774193548
$ws.Range("N3").Value = 0.64468864468864462
$ws.Range("O3").Value = 0.64642857142857135
$ws.Range("P3").Value = 0.67597765363128492
$ws.Range("Q3").Value = 0.66323024054982815
$ws.Range("R3").Value = 0.68456375838926176
$ws.Range("S3").Value = 0.67099567099567103

$ws.Range("A1:S3").Select()
